# cover-letter.docx minor rewording pass
# ---------------------------------------------------------------------
# Applies a handful of wording/formatting tweaks to the cover letter:
#  1. Title-case two words in the manuscript title; TM -> (R); "system" -> "System"
#  2. Drop "human" from "five different human populations"
#  3. Expand the "previous studies" sentence with more detail
#  4. Rewrite the "As a result..." paragraph into two clearly signposted
#     conclusions, with First/Second underlined and the key claims italicised
#     (with the key noun phrase additionally underlined)
#  5. "We look forward to hearing from you," -> "We thank you for your
#     consideration," (and relocate the _GoBack bookmark accordingly)
# ---------------------------------------------------------------------

$d = $word.ActiveDocument

function Find-Replace($range, [string]$old, [string]$new) {
    return $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

# --- 1. Title line (paragraph 5) ---------------------------------------
$p5 = $d.Paragraphs.Item(5).Range
Find-Replace $p5 "A thorough evaluation of the Languag" "A Thorough Evaluation of the Languag" | Out-Null

$p5 = $d.Paragraphs.Item(5).Range
Find-Replace $p5 "TM" "R" | Out-Null

$p5 = $d.Paragraphs.Item(5).Range
Find-Replace $p5 ") system" ") System" | Out-Null

# --- 2/3. "human" removal + expanded "previous studies" clause (paragraph 7) ---
$p7 = $d.Paragraphs.Item(7).Range
Find-Replace $p7 "five different human populations" "five different populations" | Out-Null

$p7 = $d.Paragraphs.Item(7).Range
$oldClause = "while most previous studies check only some of LENA's metrics, we check all of them."
$newClause = "while most previously published studies (including LENA" + [char]0x2019 + "s own) report accuracy and validation on only some of LENA's metrics, we check all of them."
Find-Replace $p7 $oldClause $newClause | Out-Null

# --- 4. "As a result..." paragraph rewrite (paragraph 9) ---------------
$p9 = $d.Paragraphs.Item(9).Range
$oldP9 = "As a result, we are in the best position to conclude that, at present, it appears the system is quite robust to corpora variation, and thus can be used with similar accuracy in different datasets. However, this accuracy is extremely variable across metrics."
$newP9 = "As a result, we are in a uniquely strong position to draw two key conclusions about this system. First, our results indicate that the system is quite robust to variation across corpora, and thus can be used with similar accuracy in different datasets. Second, this accuracy is extremely variable across metrics."
Find-Replace $p9 $oldP9 $newP9 | Out-Null

# Underline "First"
$r = $d.Paragraphs.Item(9).Range
$r.Find.Execute("First", $true) | Out-Null
$r.Font.Underline = 1

# Italicise "robust to variation across "
$r = $d.Paragraphs.Item(9).Range
$r.Find.Execute("robust to variation across ", $true) | Out-Null
$r.Font.Italic = 1

# Italicise + underline "corpora"
$r = $d.Paragraphs.Item(9).Range
$r.Find.Execute("corpora", $true) | Out-Null
$r.Font.Italic = 1
$r.Font.Underline = 1

# Underline "Second"
$r = $d.Paragraphs.Item(9).Range
$r.Find.Execute("Second", $true) | Out-Null
$r.Font.Underline = 1

# Italicise "accuracy is extremely variable across "
$r = $d.Paragraphs.Item(9).Range
$r.Find.Execute("accuracy is extremely variable across ", $true) | Out-Null
$r.Font.Italic = 1

# Italicise + underline "metrics"
$r = $d.Paragraphs.Item(9).Range
$r.Find.Execute("metrics", $true) | Out-Null
$r.Font.Italic = 1
$r.Font.Underline = 1

# --- 5. Closing line (paragraph 11) + _GoBack bookmark relocation ------
$p11 = $d.Paragraphs.Item(11).Range
$oldClose = [char]9 + "We look forward to hearing from you,"
$newClose = "We thank you for your consideration,"
Find-Replace $p11 $oldClose $newClose | Out-Null

$p11 = $d.Paragraphs.Item(11).Range
$p11.Find.Execute("consideration,", $true) | Out-Null
$bmPos = $p11.Start + ("consideration".Length)
$bmRange = $d.Range($bmPos, $bmPos + 1)
$bmRange.Collapse(1)
$d.Bookmarks.Add("_GoBack", $bmRange) | Out-Null

Write-Host "Done."
